$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold formatted text (e.g. "51.954.37",
# "  +0.15%  "). Excel's COM Range.Value setter auto-converts plain numeric-
# looking strings to real numbers (dropping things like trailing zeros), so
# force the whole working range to Text format before writing, then restore
# the default "Normal" style afterwards so no stray style indices are left
# behind in the saved file.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "51.954.37"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.821.42"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "355.86"
$ws.Range("E5").Value = "  +3.64%  "
$ws.Range("D6").Value = "111.70"
$ws.Range("E6").Value = "  -3.18%  "
$ws.Range("E7").Value = "  +3.75%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.600"
$ws.Range("E9").Value = "  +4.26%  "
$ws.Range("D10").Value = "40.86"
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "19.94"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "7.83"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("D15").Value = "3.261.30"
$ws.Range("D16").Value = "2.823.64"
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("D17").Value = "0.924"
$ws.Range("E17").Value = "  +4.85%  "
$ws.Range("D18").Value = "51.821.73"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  +7.69%  "
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("D21").Value = "13.41"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").Value = "0.0{0}0994" -f [char]0x2083
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("D23").Value = "70.11"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "268.05"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "2.81"
$ws.Range("E25").Value = "  +1.89%  "
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("B28").Value = "VeChain"
$ws.Range("C28").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D28").Value = "0.0513"
$ws.Range("E28").Value = "  +28.07%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "10.31"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.25"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").Value = "52.61"
$ws.Range("E32").Value = "  +4.98%  "
$ws.Range("D33").Value = "34.81"
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("D34").Value = "5.91"
$ws.Range("E34").Value = "  +3.36%  "
$ws.Range("D35").Value = "5.46"
$ws.Range("E35").Value = "  +10.49%  "
$ws.Range("E36").Value = "  +3.21%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "3.29"
$ws.Range("E38").Value = "  +2.41%  "
$ws.Range("E39").Value = "  -2.66%  "
$ws.Range("D40").Value = "18.34"
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "2.54"
$ws.Range("E42").Value = "  -4.45%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "23.26"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("D44").Value = "126.29"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").Value = "2.27"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").Value = "2.097.62"
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("D47").Value = "3.35"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("D49").Value = "6.03"
$ws.Range("E49").Value = "  +8.86%  "
$ws.Range("D50").Value = "0.987"
$ws.Range("E50").Value = "  +10.08%  "
$ws.Range("E51").Value = "  +1.71%  "

$rng.Style = "Normal"
